$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text in the Price column from Excel's automatic
# number coercion (e.g. "0.0710" -> 0.071, "1.20" -> 1.2, "1.00" -> 1) by
# forcing the data columns to Text format before writing the new values,
# then clearing that temporary formatting afterwards so the cells keep
# their original (default) style, same as in the source workbook.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '34.025.28'
$ws.Range('E2').Value = '  -1.89%  '
$ws.Range('D3').Value = '1.788.70'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '221.45'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '32.43'
$ws.Range('D9').Value = '0.284'
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').Value = '0.0710'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = '2.042.08'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '1.783.72'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '10.86'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').Value = '0.625'
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').Value = '34.017.09'
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').Value = '4.16'
$ws.Range('E17').Value = '  -3.09%  '
$ws.Range('D18').Value = '67.83'
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('D19').Value = '243.72'
$ws.Range('E19').Value = '  -3.82%  '
$ws.Range('D20').Value = '0.0₃0782'
$ws.Range('E20').Value = '  -2.24%  '
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '10.83'
$ws.Range('E22').Value = '  +1.63%  '
$ws.Range('D23').Value = '4.08'
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('D25').Value = '157.69'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').Value = '16.34'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '7.04'
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('E28').Value = '  -1.53%  '
$ws.Range('E29').Value = '  -0.10%  '
$ws.Range('D30').Value = '0.0519'
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('E33').Value = '  -2.75%  '
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('D35').Value = '1.395.40'
$ws.Range('E35').Value = '  -2.83%  '
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('E38').Value = '  -3.25%  '
$ws.Range('D39').Value = '79.63'
$ws.Range('E39').Value = '  -5.85%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '0.921'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').Value = '2.35'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('E42').Value = '  -2.57%  '
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('D45').Value = '107.29'
$ws.Range('E45').Value = '  +2.02%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '5.88'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').Value = '0.0492'
$ws.Range('E47').Value = '  +0.46%  '
$ws.Range('D48').Value = '1.944.82'
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').Value = '11.94'
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').Value = '0.0₆0127'
$ws.Range('E51').Value = '  +0.99%  '

# Clear the temporary Text-format styling (values keep their text type).
$dataRange.ClearFormats()
